$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 322; this pushes the existing rows
# 322..413 down to 323..414 (dimension grows from R413 to R414).
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row 322 with its data.
$ws.Range("A322").Value = 11
$ws.Range("B322").Value = "Vega Monumental Concepción"
$ws.Range("C322").Value = "Bíobío"
$ws.Range("D322").Value = 44932
$ws.Range("E322").Value = 8
$ws.Range("F322").Value = 100112023
$ws.Range("G322").Value = "Brócoli"
$ws.Range("H322").Value = "Sin especificar"
$ws.Range("I322").Value = "Primera"
$ws.Range("J322").Value = 2000
$ws.Range("K322").Value = 650
$ws.Range("L322").Value = 700
$ws.Range("M322").Value = 675
$ws.Range("N322").Value = "$/unidad"
$ws.Range("O322").Value = "Región Metropolitana"
$ws.Range("P322").Value = 675
$ws.Range("Q322").Value = 1
$ws.Range("R322").Value = "Hortaliza"

# Match the date-cell number format used by the other rows in column D.
$ws.Range("D322").NumberFormat = $ws.Range("D323").NumberFormat
